$d = $word.ActiveDocument
$d.Content.Find.Execute("Pull ", $true, $false, $false, $false, $false, $true, 1, $false, "Get ", 2)
